$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (item id 5512)
$ws.Range("H33").Value = 1994.4286
$ws.Range("J33").Value = 380.33334
$ws.Range("L33").Value = 380.33334
$ws.Range("N33").Value = -838.33334

# Row 132 (item id 44049)
$ws.Range("H132").Value = 2101.4854
$ws.Range("I132").Value = 1830.6102
$ws.Range("J132").Value = 3877.2222
$ws.Range("K132").Value = 5491.8306
$ws.Range("L132").Value = 11631.6666
$ws.Range("M132").Value = -2961.8306
$ws.Range("N132").Value = -16691.6666

# Row 137 (item id 44013)
$ws.Range("H137").Value = 840183.1
$ws.Range("I137").Value = 3248.8333
$ws.Range("J137").Value = 3350986
$ws.Range("K137").Value = 9746.499899999999
$ws.Range("L137").Value = 10052958
$ws.Range("M137").Value = -7196.499899999999
$ws.Range("N137").Value = -10058058

# Row 138 (item id 44169)
$ws.Range("H138").Value = 2262.5
$ws.Range("I138").Value = 1098.5625
$ws.Range("J138").Value = 3955.5
$ws.Range("K138").Value = 3295.6875
$ws.Range("L138").Value = 11866.5
$ws.Range("M138").Value = 1844.3125
$ws.Range("N138").Value = -22146.5

# Row 141 (item id 44161)
$ws.Range("H141").Value = 1374.5264
$ws.Range("I141").Value = 1339.7778
$ws.Range("K141").Value = 4019.3334
$ws.Range("M141").Value = 1160.6666

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (item id 44147)
$ws.Range("H32").Value = 5348.52
$ws.Range("I32").Value = 3542.4106
$ws.Range("K32").Value = 3542.4106
$ws.Range("M32").Value = -3255.4106

# Row 61 (item id 43999)
$ws.Range("H61").Value = 23790.156
$ws.Range("I61").Value = 2753.3547
$ws.Range("J61").Value = 70371.64
$ws.Range("K61").Value = 2753.3547
$ws.Range("L61").Value = 70371.64
$ws.Range("M61").Value = -2541.3547
$ws.Range("N61").Value = -70795.64

# Row 74 (item id 44000)
$ws.Range("H74").Value = 181628.9
$ws.Range("I74").Value = 114036.25
$ws.Range("J74").Value = 451999.5
$ws.Range("K74").Value = 114036.25
$ws.Range("L74").Value = 451999.5
$ws.Range("M74").Value = -113162.25
$ws.Range("N74").Value = -453747.5

# Row 77 (item id 44000)
$ws.Range("H77").Value = 181628.9
$ws.Range("I77").Value = 114036.25
$ws.Range("J77").Value = 451999.5
$ws.Range("K77").Value = 570181.25
$ws.Range("L77").Value = 2259997.5
$ws.Range("M77").Value = -565813.25
$ws.Range("N77").Value = -2268733.5

# Row 123 (item id 34107)
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 132 (item id 43997)
$ws.Range("H132").Value = 3083.0322
$ws.Range("I132").Value = 2490.8076
$ws.Range("J132").Value = 6162.6
$ws.Range("K132").Value = 7472.4228
$ws.Range("L132").Value = 18487.8
$ws.Range("M132").Value = -4942.4228
$ws.Range("N132").Value = -23547.8

# Row 136 (item id 43999)
$ws.Range("H136").Value = 23790.156
$ws.Range("I136").Value = 2753.3547
$ws.Range("J136").Value = 70371.64
$ws.Range("K136").Value = 8260.0641
$ws.Range("L136").Value = 211114.92
$ws.Range("M136").Value = -5710.0641
$ws.Range("N136").Value = -216214.92

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (item id 44023)
$ws.Range("H31").Value = 2972.8135
$ws.Range("I31").Value = 2664.3572
$ws.Range("J31").Value = 3068.7778
$ws.Range("K31").Value = 2664.3572
$ws.Range("L31").Value = 3068.7778
$ws.Range("M31").Value = -2369.3572
$ws.Range("N31").Value = -3658.7778

# Row 34 (item id 44023)
$ws.Range("H34").Value = 2972.8135
$ws.Range("I34").Value = 2664.3572
$ws.Range("J34").Value = 3068.7778
$ws.Range("K34").Value = 2664.3572
$ws.Range("L34").Value = 3068.7778
$ws.Range("M34").Value = -2462.3572
$ws.Range("N34").Value = -3472.7778

# Row 58 (item id 44021)
$ws.Range("H58").Value = 2775.275
$ws.Range("I58").Value = 2551.4849
$ws.Range("J58").Value = 3830.2856
$ws.Range("K58").Value = 2551.4849
$ws.Range("L58").Value = 3830.2856
$ws.Range("M58").Value = -2348.4849
$ws.Range("N58").Value = -4236.2856

# Row 132 (item id 44019)
$ws.Range("H132").Value = 4740.303
$ws.Range("I132").Value = 1627.5927
$ws.Range("J132").Value = 18747.5
$ws.Range("K132").Value = 4882.7781
$ws.Range("L132").Value = 56242.5
$ws.Range("M132").Value = -2352.7781
$ws.Range("N132").Value = -61302.5

# Row 134 (item id 44020)
$ws.Range("H134").Value = 1825.4912
$ws.Range("I134").Value = 1612.0426
$ws.Range("J134").Value = 2828.7
$ws.Range("K134").Value = 4836.1278
$ws.Range("L134").Value = 8486.099999999999
$ws.Range("M134").Value = -2301.1278
$ws.Range("N134").Value = -13556.1

# Row 136 (item id 44021)
$ws.Range("H136").Value = 2775.275
$ws.Range("I136").Value = 2551.4849
$ws.Range("J136").Value = 3830.2856
$ws.Range("K136").Value = 7654.4547
$ws.Range("L136").Value = 11490.8568
$ws.Range("M136").Value = -5104.4547
$ws.Range("N136").Value = -16590.8568

$ws = $wb.Worksheets.Item("CUL")
# Row 122 (item id 36078)
$ws.Range("H122").Value = 41872.047
$ws.Range("I122").Value = 637.7143
$ws.Range("J122").Value = 61114.734
$ws.Range("K122").Value = 5739.428699999999
$ws.Range("L122").Value = 550032.6059999999
$ws.Range("M122").Value = -3289.428699999999
$ws.Range("N122").Value = -554932.6059999999

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (item id 27710)
$ws.Range("H113").Value = 5138.1665
$ws.Range("I113").Value = 5303.3687
$ws.Range("K113").Value = 5303.3687
$ws.Range("M113").Value = -3133.3687

# Row 132 (item id 44008)
$ws.Range("H132").Value = 3598.6619
$ws.Range("I132").Value = 3595.623
$ws.Range("J132").Value = 3617.2
$ws.Range("K132").Value = 10786.869
$ws.Range("L132").Value = 10851.6
$ws.Range("M132").Value = -8256.869000000001
$ws.Range("N132").Value = -15911.6

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (item id 36248)
$ws.Range("H40").Value = 3953.575
$ws.Range("I40").Value = 3521.1177
$ws.Range("J40").Value = 6404.1665
$ws.Range("K40").Value = 3521.1177
$ws.Range("L40").Value = 6404.1665
$ws.Range("M40").Value = -3385.1177
$ws.Range("N40").Value = -6676.1665

# Row 132 (item id 44058)
$ws.Range("H132").Value = 2616.9644
$ws.Range("I132").Value = 2191.34
$ws.Range("J132").Value = 6163.8335
$ws.Range("K132").Value = 6574.02
$ws.Range("L132").Value = 18491.5005
$ws.Range("M132").Value = -4044.02
$ws.Range("N132").Value = -23551.5005

# Row 136 (item id 44060)
$ws.Range("H136").Value = 3245.7778
$ws.Range("I136").Value = 2725.875
$ws.Range("J136").Value = 7405
$ws.Range("K136").Value = 8177.625
$ws.Range("L136").Value = 22215
$ws.Range("M136").Value = -5627.625
$ws.Range("N136").Value = -27315

$ws = $wb.Worksheets.Item("WVR")
# Row 136 (item id 44031)
$ws.Range("H136").Value = 3491.1929
$ws.Range("I136").Value = 3253.157
$ws.Range("J136").Value = 5514.5
$ws.Range("K136").Value = 9759.471000000001
$ws.Range("L136").Value = 16543.5
$ws.Range("M136").Value = -7209.471000000001
$ws.Range("N136").Value = -21643.5
